$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt2"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01070233333333333
$ws.Range("H2").Value = 0.032107
$ws.Range("I2").Value = 0.006017198313602724
$ws.Range("J2").Value = 0.006017198313602724
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 27.25040933333333
$ws.Range("N2").Value = 81.751228
$ws.Range("O2").Value = 0.455635031912059
$ws.Range("P2").Value = 0.4556350319120589
$ws.Range("Q2").Value = 0.2916429641551111
$ws.Range("R2").Value = 2.624786677396
$ws.Range("S2").Value = 0.002741646345639565
$ws.Range("T2").Value = 0.002741646345639564

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt2"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01070233333333333
$ws.Range("H3").Value = 0.032107
$ws.Range("I3").Value = 0.006017198313602724
$ws.Range("J3").Value = 0.006017198313602724
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.03147833333333
$ws.Range("N3").Value = 42.094435
$ws.Range("O3").Value = 0.2346105337346748
$ws.Range("P3").Value = 0.2346105337346748
$ws.Range("Q3").Value = 0.1501695582827778
$ws.Range("R3").Value = 1.351526024545
$ws.Range("S3").Value = 0.001411698107941721
$ws.Range("T3").Value = 0.00141169810794172

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01070233333333333
$ws.Range("H4").Value = 0.032107
$ws.Range("I4").Value = 0.006017198313602724
$ws.Range("J4").Value = 0.006017198313602724
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3744183333333334
$ws.Range("N4").Value = 1.123255
$ws.Range("O4").Value = 0.006260387033823881
$ws.Range("P4").Value = 0.00626038703382388
$ws.Range("Q4").Value = 0.004007149809444444
$ws.Range("R4").Value = 0.036064348285
$ws.Range("S4").Value = [double]"3.766999030242542e-05"
$ws.Range("T4").Value = [double]"3.766999030242541e-05"

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01070233333333333
$ws.Range("H5").Value = 0.032107
$ws.Range("I5").Value = 0.006017198313602724
$ws.Range("J5").Value = 0.006017198313602724
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1421396666666667
$ws.Range("N5").Value = 0.426419
$ws.Range("O5").Value = 0.002376617934997971
$ws.Range("P5").Value = 0.00237661793499797
$ws.Range("Q5").Value = 0.001521226092555556
$ws.Range("R5").Value = 0.013691034833
$ws.Range("S5").Value = [double]"1.430058143054778e-05"
$ws.Range("T5").Value = [double]"1.430058143054778e-05"

$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01070233333333333
$ws.Range("H6").Value = 0.032107
$ws.Range("I6").Value = 0.006017198313602724
$ws.Range("J6").Value = 0.006017198313602724
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 18.009092
$ws.Range("N6").Value = 54.027276
$ws.Range("O6").Value = 0.3011174293844444
$ws.Range("P6").Value = 0.3011174293844444
$ws.Range("Q6").Value = 0.1927393056146667
$ws.Range("R6").Value = 1.734653750532
$ws.Range("S6").Value = 0.001811883288288466
$ws.Range("T6").Value = 0.001811883288288466

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt2"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.767921666666667
$ws.Range("H7").Value = 5.303765
$ws.Range("I7").Value = 0.9939828016863973
$ws.Range("J7").Value = 0.9939828016863973
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.25040933333333
$ws.Range("N7").Value = 81.751228
$ws.Range("O7").Value = 0.455635031912059
$ws.Range("P7").Value = 0.4556350319120589
$ws.Range("Q7").Value = 48.17658908593556
$ws.Range("R7").Value = 433.58930177342
$ws.Range("S7").Value = 0.4528933855664194
$ws.Range("T7").Value = 0.4528933855664193

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt2"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.767921666666667
$ws.Range("H8").Value = 5.303765
$ws.Range("I8").Value = 0.9939828016863973
$ws.Range("J8").Value = 0.9939828016863973
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.03147833333333
$ws.Range("N8").Value = 42.094435
$ws.Range("O8").Value = 0.2346105337346748
$ws.Range("P8").Value = 0.2346105337346748
$ws.Range("Q8").Value = 24.80655456086389
$ws.Range("R8").Value = 223.258991047775
$ws.Range("S8").Value = 0.2331988356267331
$ws.Range("T8").Value = 0.2331988356267331

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt2"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.767921666666667
$ws.Range("H9").Value = 5.303765
$ws.Range("I9").Value = 0.9939828016863973
$ws.Range("J9").Value = 0.9939828016863973
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3744183333333334
$ws.Range("N9").Value = 1.123255
$ws.Range("O9").Value = 0.006260387033823881
$ws.Range("P9").Value = 0.00626038703382388
$ws.Range("Q9").Value = 0.6619422838972223
$ws.Range("R9").Value = 5.957480555075001
$ws.Range("S9").Value = 0.006222717043521456
$ws.Range("T9").Value = 0.006222717043521456

$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt2"
$ws.Range("C10").Value = "Fzd4"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.767921666666667
$ws.Range("H10").Value = 5.303765
$ws.Range("I10").Value = 0.9939828016863973
$ws.Range("J10").Value = 0.9939828016863973
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1421396666666667
$ws.Range("N10").Value = 0.426419
$ws.Range("O10").Value = 0.002376617934997971
$ws.Range("P10").Value = 0.00237661793499797
$ws.Range("Q10").Value = 0.2512917963927778
$ws.Range("R10").Value = 2.261626167535
$ws.Range("S10").Value = 0.002362317353567423
$ws.Range("T10").Value = 0.002362317353567422

$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt2"
$ws.Range("C11").Value = "Fzd4"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.767921666666667
$ws.Range("H11").Value = 5.303765
$ws.Range("I11").Value = 0.9939828016863973
$ws.Range("J11").Value = 0.9939828016863973
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 18.009092
$ws.Range("N11").Value = 54.027276
$ws.Range("O11").Value = 0.3011174293844444
$ws.Range("P11").Value = 0.3011174293844444
$ws.Range("Q11").Value = 31.83866394379334
$ws.Range("R11").Value = 286.54797549414
$ws.Range("S11").Value = 0.2993055460961559
$ws.Range("T11").Value = 0.2993055460961559

